# Add a new "Sexy Header" metadata row to the "#Quotes" sheet (carpenter's
# list of tables scripted by sexy now also includes the quotes table's
# sexy header), and switch the active/selected tab from "#Quotes" to
# "Quotes Table".

$wb = $excel.ActiveWorkbook

# --- "#Quotes" sheet: append the new metadata row -------------------------
$quotes = $wb.Worksheets.Item("#Quotes")

$quotes.Range("A11").Value = "Sexy Header"
$quotes.Range("B11").Value = "tables\rococo.tables.test.sxh"

# Move this sheet's own selection off of the old B8 cell onto C1, and drop
# its "active tab" status in favour of "Quotes Table" below.
$quotes.Range("C1").Select()

# --- "Quotes Table" sheet: becomes the active / selected tab --------------
$quotesTable = $wb.Worksheets.Item("Quotes Table")
$quotesTable.Activate()
